# "add 'import csv file' function and edit save, load student functions"
#
# Concretely (per the OOXML diff) this commit:
#   1. Renames the sheet tab from "Sheet1" to "19CLC9 Students".
#   2. Re-orders the Gender / DoB / Class columns so Gender (D) comes right
#      after Fullname, followed by DoB (E) and Class (F).
#   3. Gives the (new) Gender column D a centered horizontal alignment.
#   4. Removes the trailing blank rows (7-11) that used to pad the sheet.
#   5. Leaves the active selection on F11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename the worksheet tab.
# ---------------------------------------------------------------------
$ws.Name = "19CLC9 Students"

# ---------------------------------------------------------------------
# 2. Drop the leftover empty rows (7:11) under the real data (rows 1:6).
# ---------------------------------------------------------------------
$ws.Range("A7:F11").Delete() | Out-Null

# ---------------------------------------------------------------------
# 3. Re-order columns D (DoB) / E (Class) / F (Gender) into
#    D (Gender) / E (DoB) / F (Class).
#
#    Move the DoB number-format (applied to D2:D6) over to the new DoB
#    column (E2:E6) first, via a scratch column, so the style survives
#    the swap intact, then swap the raw values column-by-column.
# ---------------------------------------------------------------------
$ws.Range("D2:D6").Copy() | Out-Null
$ws.Range("H2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (keep the DoB format aside)

for ($r = 1; $r -le 6; $r++) {
    $dCell = $ws.Cells.Item($r, 4)   # DoB
    $eCell = $ws.Cells.Item($r, 5)   # Class
    $fCell = $ws.Cells.Item($r, 6)   # Gender

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    $fCell.Value2 = $eVal   # Class moves to F
    $eCell.Value2 = $dVal   # DoB moves to E
    $dCell.Value2 = $fVal   # Gender moves to D
}

# Column D used to hold DoB values (date number format) - reset it to the
# default format before giving it its own (centered) look.
$ws.Range("D1:D6").ClearFormats() | Out-Null

# Re-apply the DoB number format (captured above) onto the new DoB column E.
$ws.Range("H2:H6").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H2:H6").Clear() | Out-Null

# ---------------------------------------------------------------------
# 4. Center-align the new Gender column (D).
# ---------------------------------------------------------------------
$ws.Range("D1:D6").HorizontalAlignment = -4108   # xlCenter

# ---------------------------------------------------------------------
# 5. Leave the selection on F11, matching the saved view state.
# ---------------------------------------------------------------------
$ws.Range("F11").Select() | Out-Null
